$wb = $excel.ActiveWorkbook

# --- clients sheet: scroll the viewport down (best effort; selection stays at F6) ---
$clients = $wb.Worksheets.Item("clients")
$clients.Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$clients.Range("F6").Select() | Out-Null

# --- vendor_inventory sheet: new trailing (empty) row + selection change ---
$vendorInv = $wb.Worksheets.Item("vendor_inventory")
$vendorInv.Select() | Out-Null
$vendorInv.Range("A14").Style = "Normal"
$vendorInv.Range("C14").Select() | Out-Null

# --- expense_reports sheet: the "Items Out of Stock"/"Items Not Found" columns moved out ---
$expenseReports = $wb.Worksheets.Item("expense_reports")
$expenseReports.Select() | Out-Null
$expenseReports.Range("F1").ClearContents() | Out-Null
$expenseReports.Range("G1").ClearContents() | Out-Null
$expenseReports.Range("D12").Select() | Out-Null

# --- add the new out_of_stock sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$outOfStock = $wb.Worksheets.Add($null, $lastSheet)
$outOfStock.Name = "out_of_stock"
$outOfStock.Range("A1").Value = "Vendor"
$outOfStock.Range("B1").Value = "Items Out of Stock"
$outOfStock.Range("C1").Value = "Items Not Found"
$outOfStock.Range("D13").Select() | Out-Null
